# Actualización automática 2025-12-07 17:44:30
# Applies updated sales / budget figures to the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (product group sales) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# DECORHOME S.C.C. now shows a PORCELANATO sale of 1432.08
$wsGrupo.Range("M9").Value = 1432.08

# Summary row: count of clients with sales in PORCELANATO column goes from 2 to 3
$wsGrupo.Range("M24").Value = "3 de 22"

# --- Sheet "VENTA MENSUAL" (monthly sales / budget) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# PRESUPUESTO (budget) column updates
$wsMensual.Range("G3").Value = 0
$wsMensual.Range("G5").Value = 5000
$wsMensual.Range("G7").Value = 5000
$wsMensual.Range("G8").Value = 500

# DECORHOME S.C.C. now has a December sale and a revised budget
$wsMensual.Range("F9").Value = 1484.26
$wsMensual.Range("G9").Value = 5000

$wsMensual.Range("G11").Value = 1000
$wsMensual.Range("G12").Value = 6000
$wsMensual.Range("G14").Value = 6500
$wsMensual.Range("G15").Value = 6500
$wsMensual.Range("G16").Value = 0
$wsMensual.Range("G17").Value = 1500
$wsMensual.Range("G18").Value = 5000
$wsMensual.Range("G21").Value = 0
$wsMensual.Range("G22").Value = 0
$wsMensual.Range("G23").Value = 7000

# Totals row (recomputed sums)
$wsMensual.Range("F24").Value = 11949.31
$wsMensual.Range("G24").Value = 49000
